$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows to append (rows 83-86)
$data = @(
    @(44314, 485,  3566, 1300, 5601, 131, 1666, 1000, 28500, 0, 0, 0,    3638, 144, 1081, 42971),
    @(44315, 1270, 3616, 800,  5651, 188, 1835, 0,    28500, 0, 0, 0,    3638, 597, 1138, 43241),
    @(44316, 40,   3531, 850,  5601, 725, 1946, 1000, 28500, 0, 0, 0,    3638, 11,  1126, 43217),
    @(44319, 433,  3333, 740,  5751, 683, 1929, 1000, 28500, 0, 0, 0,    3638, 182, 1124, 43151)
)

$startRow = 83
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
